# Insert two new daily-price rows for Acelga (date serial 45021) into the
# dataset. Excel shifts every existing row from 1087 onward down by two rows
# when the new rows are inserted, which is the bulk of the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows right before the current row 1087.
$ws.Range("A1087:A1088").EntireRow.Insert()

# New row 1087: Acelga, Primera, Región Metropolitana
$ws.Range("A1087").Value = 6
$ws.Range("B1087").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1087").Value = "Metropolitana"
$ws.Range("D1087").Value = 45021
$ws.Range("E1087").Value = 13
$ws.Range("F1087").Value = 100112009
$ws.Range("G1087").Value = "Acelga"
$ws.Range("H1087").Value = "Sin especificar"
$ws.Range("I1087").Value = "Primera"
$ws.Range("J1087").Value = 280
$ws.Range("K1087").Value = 12000
$ws.Range("L1087").Value = 12000
$ws.Range("M1087").Value = 12000
$ws.Range("N1087").Value = "`$/docena de atados"
$ws.Range("O1087").Value = "Región Metropolitana"
$ws.Range("P1087").Value = 4000
$ws.Range("Q1087").Value = 3
$ws.Range("R1087").Value = "Hortaliza"

# New row 1088: Acelga, Segunda, Región Metropolitana
$ws.Range("A1088").Value = 6
$ws.Range("B1088").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1088").Value = "Metropolitana"
$ws.Range("D1088").Value = 45021
$ws.Range("E1088").Value = 13
$ws.Range("F1088").Value = 100112009
$ws.Range("G1088").Value = "Acelga"
$ws.Range("H1088").Value = "Sin especificar"
$ws.Range("I1088").Value = "Segunda"
$ws.Range("J1088").Value = 170
$ws.Range("K1088").Value = 9000
$ws.Range("L1088").Value = 9000
$ws.Range("M1088").Value = 9000
$ws.Range("N1088").Value = "`$/docena de atados"
$ws.Range("O1088").Value = "Región Metropolitana"
$ws.Range("P1088").Value = 3000
$ws.Range("Q1088").Value = 3
$ws.Range("R1088").Value = "Hortaliza"
